$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", formatted like the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column H (rows 2-12): 0/1 flag values
$values = @(1, 0, 1, 1, 0, 1, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}

$excel.CutCopyMode = 0
